$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.782.63"
$ws.Range("E2").Value = "  +4.26%  "
$ws.Range("D3").Value = "2.251.63"
$ws.Range("E3").Value = "  +3.66%  "
$ws.Range("E4").Value = "  +0.09%  "
$__style = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "249.92"
$ws.Range("D5").Style = $__style
$ws.Range("E5").Value = "  +0.47%  "
$__style = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.624"
$ws.Range("D6").Style = $__style
$ws.Range("E6").Value = "  +0.93%  "
$__style = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "70.61"
$ws.Range("D7").Style = $__style
$ws.Range("E7").Value = "  +5.40%  "
$ws.Range("B8").Value = "Cardano"
$ws.Range("C8").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$__style = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.671"
$ws.Range("D8").Style = $__style
$ws.Range("E8").Value = "  +18.62%  "
$ws.Range("B9").Value = "USDC"
$ws.Range("C9").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$__style = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.00"
$ws.Range("D9").Style = $__style
$ws.Range("E9").Value = "  -0.08%  "
$__style = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.57"
$ws.Range("D10").Style = $__style
$ws.Range("E10").Value = "  +11.47%  "
$ws.Range("E11").Value = "  +4.76%  "
$__style = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "59.16"
$ws.Range("D12").Style = $__style
$ws.Range("E12").Value = "  +1.61%  "
$__style = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.54"
$ws.Range("D13").Style = $__style
$ws.Range("E13").Value = "  +9.87%  "
$__style = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.105"
$ws.Range("D14").Style = $__style
$ws.Range("E14").Value = "  +1.45%  "
$ws.Range("D15").Value = "2.588.21"
$ws.Range("E15").Value = "  +3.67%  "
$ws.Range("B16").Value = "Polygon"
$ws.Range("C16").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$__style = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.883"
$ws.Range("D16").Style = $__style
$ws.Range("E16").Value = "  +3.21%  "
$ws.Range("B17").Value = "Chainlink"
$ws.Range("C17").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$__style = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.76"
$ws.Range("D17").Style = $__style
$ws.Range("E17").Value = "  +4.29%  "
$ws.Range("D18").Value = "2.237.24"
$ws.Range("E18").Value = "  +2.13%  "
$ws.Range("D19").Value = "42.694.77"
$ws.Range("E19").Value = "  +4.31%  "
$ws.Range("D20").Value = "0.0₃0989"
$ws.Range("E20").Value = "  +5.57%  "
$__style = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.28"
$ws.Range("D21").Style = $__style
$ws.Range("E21").Value = "  +3.37%  "
$__style = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "72.96"
$ws.Range("D22").Style = $__style
$ws.Range("E22").Value = "  +2.15%  "
$__style = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "233.48"
$ws.Range("D23").Style = $__style
$ws.Range("E23").Value = "  +1.56%  "
$__style = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.07"
$ws.Range("D24").Style = $__style
$ws.Range("E24").Value = "  -0.17%  "
$__style = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.94"
$ws.Range("D25").Style = $__style
$ws.Range("E25").Value = "  +5.07%  "
$ws.Range("B26").Value = "Dai"
$ws.Range("C26").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$__style = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").Style = $__style
$ws.Range("E26").Value = "  +0.08%  "
$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$__style = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.36"
$ws.Range("D27").Style = $__style
$ws.Range("E27").Value = "  +0.16%  "
$ws.Range("E28").Value = "  +0.79%  "
$ws.Range("E29").Value = "  -1.67%  "
$ws.Range("E30").Value = "  +10.01%  "
$__style = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "167.47"
$ws.Range("D31").Style = $__style
$ws.Range("E31").Value = "  -0.26%  "
$ws.Range("E32").Value = "  +3.70%  "
$__style = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.49"
$ws.Range("D33").Style = $__style
$ws.Range("E33").Value = "  +14.33%  "
$__style = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.128"
$ws.Range("D34").Style = $__style
$ws.Range("E34").Value = "  +7.69%  "
$__style = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0789"
$ws.Range("D35").Style = $__style
$ws.Range("E35").Value = "  +7.18%  "
$__style = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "30.93"
$ws.Range("D36").Style = $__style
$ws.Range("E36").Value = "  +21.04%  "
$ws.Range("E37").Value = "  +4.47%  "
$__style = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.40"
$ws.Range("D38").Style = $__style
$ws.Range("E38").Value = "  +8.12%  "
$ws.Range("E39").Value = "  +3.46%  "
$ws.Range("E40").Value = "  +7.15%  "
$__style = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.32"
$ws.Range("D41").Style = $__style
$ws.Range("E41").Value = "  +7.24%  "
$__style = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "12.55"
$ws.Range("D42").Style = $__style
$ws.Range("E42").Value = "  +10.36%  "
$ws.Range("E43").Value = "  +6.31%  "
$__style = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "62.73"
$ws.Range("D44").Style = $__style
$ws.Range("E44").Value = "  +2.48%  "
$__style = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.201"
$ws.Range("D45").Style = $__style
$ws.Range("E45").Value = "  +5.21%  "
$ws.Range("E46").Value = "  +5.05%  "
$__style = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.85"
$ws.Range("D47").Style = $__style
$ws.Range("E47").Value = "  +1.46%  "
$ws.Range("E48").Value = "  +4.98%  "
$ws.Range("E49").Value = "  -0.18%  "
$__style = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.18"
$ws.Range("D50").Style = $__style
$ws.Range("E50").Value = "  +0.85%  "
$ws.Range("E51").Value = "  +4.35%  "
